$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Title (appears twice: H1 heading and bold paragraph near the end)
Replace-Text "Play Hit it Hard Free: Review of ELK's Fruit Machine-inspired Slot" "Play Hit it Hard for Free - Exciting Vintage-Style Slot Game"

# "What we like" bullet list
Replace-Text "Great vintage-style graphics and design" "Unique blend of vintage style and modern elements"
Replace-Text "Various bonus features like multipliers, free spins, and wild spins" "Decent payouts and winning possibilities"
Replace-Text "Flexible betting range from 20 cents to €100 per spin" "Various bonus features including multipliers and free spins"
Replace-Text "Decent payouts with a maximum payout of 2,500 times the bet amount" "Well-executed mix of modern and vintage styling"

# "What we don't like" bullet list
Replace-Text "Only five pay lines may not be enough for some players" "Limited number of pay lines"
Replace-Text "The game may feel too traditional for those seeking modern features" "Relatively low maximum payout"

# Meta description (italic paragraph)
Replace-Text "Experience the nostalgia of classic Fruit Machines with Hit it Hard, the online slot game by ELK. Free play review of payouts, features, and design." "Read our review of Hit it Hard and play this vintage-style slot game for free. Enjoy bonus features and decent payouts."

Write-Output "Replacements complete"
